# Practice with Excel and Data Tables
# Fill in the "Total" (Quantity * Price) and "Big order" columns on the
# Orders sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Row 2: OrderID 1, Oranges - Qty 7 x Price 2 = 14 (big order)
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = "Yes"

# Row 3: OrderID 2, Apples - Qty 6 x Price 1 = 6
$ws.Range("E3").Value = 6

# Row 4: OrderID 3, Bananas - Qty 1 x Price 1 = 1
$ws.Range("E4").Value = 1

# Row 5: OrderID 4, Peaches - Qty 4 x Price 3 = 12 (big order)
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = "Yes"

# Leave the selection where the author last clicked
$ws.Range("J19").Select()
